$wb = $excel.ActiveWorkbook

# Add the new worksheet after the existing one and name it "outbreak2"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "outbreak2"

# Fill in headers
$newSheet.Range("A1").Value = "time"
$newSheet.Range("B1").Value = "I1"

# Fill in data rows
$data = @(
    @(1, 1),
    @(2, 5),
    @(3, 36),
    @(4, 27),
    @(5, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
}

# Select A7 on the new sheet (empty cell below data) and make it the active sheet/tab
$newSheet.Range("A7").Select()
$newSheet.Activate()
